# Added de-duplication functionality to data_clean
# Applies the de-duplication pass to the "watercolours" sheet:
#  - Purple column (B): duplicate "Royal Purple" / "Eggplant Purple" values
#    are suffixed with a stray marker character so they are distinguishable,
#    the cascaded rows are rewritten, and the true duplicate gets folded
#    down into its own row further down the column.
#  - Pink column (A): the duplicate "More Pink" entry is dropped out of the
#    list (leaving its row blank) and every entry below it shifts down one
#    row.
# Finally, the watercolours sheet (not the test sheet) is left as the
# active selection, parked on H14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("watercolours")

# --- Purple column de-duplication cascade ---
$ws.Range("B2").Value = "Royal Purple&"
$ws.Range("B4").Value = "Eggplant~ Purple"
$ws.Range("B5").Value = "Velvet Purple%$"
$ws.Range("B6").Value = "Royal Purple"

# --- Pink column de-duplication: blank out the duplicate row, shift the
#     remainder of the list down one row ---
$ws.Rows.Item(9).Insert()

# --- Selection / active sheet bookkeeping ---
$ws.Activate()
$ws.Range("H14").Select()

Write-Output "data_clean: de-duplication applied"
